# Update the cryptocurrency price-table rows (Coin/Link/Price/Volume(1h))
# to reflect the latest scrape. Each entry below is a literal cell ->
# new-text-value pair taken from the source data feed; Price (column D)
# and Volume (column E) values are free-text strings (not numbers), so
# numeric-looking Price values are written with a Text number format to
# stop Excel's COM layer from auto-coercing them (e.g. "97.01" -> 97.01,
# or "0.0360" -> 0.036, losing the trailing zero / string type). The
# Text format nudge is cleared again immediately after the write so the
# cell's formatting stays exactly as it was before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '43.749.72' },
    @{ Cell = 'E2'; Value = '  -0.08%  ' },
    @{ Cell = 'D3'; Value = '2.318.46' },
    @{ Cell = 'E3'; Value = '  +3.66%  ' },
    @{ Cell = 'E4'; Value = '  +0.07%  ' },
    @{ Cell = 'D5'; Value = '97.01' },
    @{ Cell = 'E5'; Value = '  +4.71%  ' },
    @{ Cell = 'D6'; Value = '271.75' },
    @{ Cell = 'E6'; Value = '  +0.12%  ' },
    @{ Cell = 'E7'; Value = '  +0.99%  ' },
    @{ Cell = 'E8'; Value = '  +0.01%  ' },
    @{ Cell = 'D9'; Value = '0.624' },
    @{ Cell = 'E9'; Value = '  +0.05%  ' },
    @{ Cell = 'D10'; Value = '45.38' },
    @{ Cell = 'E10'; Value = '  -2.21%  ' },
    @{ Cell = 'E11'; Value = '  -1.53%  ' },
    @{ Cell = 'D12'; Value = '8.05' },
    @{ Cell = 'E12'; Value = '  -2.72%  ' },
    @{ Cell = 'D13'; Value = '0.106' },
    @{ Cell = 'E13'; Value = '  +0.43%  ' },
    @{ Cell = 'D14'; Value = '2.655.84' },
    @{ Cell = 'E14'; Value = '  +3.27%  ' },
    @{ Cell = 'E15'; Value = '  +2.91%  ' },
    @{ Cell = 'D16'; Value = '0.878' },
    @{ Cell = 'E16'; Value = '  +9.39%  ' },
    @{ Cell = 'D17'; Value = '2.292.64' },
    @{ Cell = 'E17'; Value = '  +2.09%  ' },
    @{ Cell = 'D18'; Value = '43.706.81' },
    @{ Cell = 'E18'; Value = '  -0.17%  ' },
    @{ Cell = 'E19'; Value = '  +4.46%  ' },
    @{ Cell = 'D20'; Value = '6.42' },
    @{ Cell = 'E20'; Value = '  +5.45%  ' },
    @{ Cell = 'D21'; Value = '73.39' },
    @{ Cell = 'E21'; Value = '  +3.66%  ' },
    @{ Cell = 'D22'; Value = '240.63' },
    @{ Cell = 'E22'; Value = '  +2.75%  ' },
    @{ Cell = 'D23'; Value = '2.26' },
    @{ Cell = 'E23'; Value = '  -3.43%  ' },
    @{ Cell = 'D24'; Value = '9.44' },
    @{ Cell = 'E24'; Value = '  +3.41%  ' },
    @{ Cell = 'E25'; Value = '  -0.08%  ' },
    @{ Cell = 'E26'; Value = '  +0.92%  ' },
    @{ Cell = 'E27'; Value = '  -0.53%  ' },
    @{ Cell = 'E28'; Value = '  -1.19%  ' },
    @{ Cell = 'D29'; Value = '2.28' },
    @{ Cell = 'E29'; Value = '  +0.68%  ' },
    @{ Cell = 'D30'; Value = '38.17' },
    @{ Cell = 'E30'; Value = '  -8.42%  ' },
    @{ Cell = 'D31'; Value = '22.43' },
    @{ Cell = 'E31'; Value = '  +6.91%  ' },
    @{ Cell = 'D32'; Value = '175.13' },
    @{ Cell = 'E32'; Value = '  +1.42%  ' },
    @{ Cell = 'D33'; Value = '0.0914' },
    @{ Cell = 'E33'; Value = '  -0.95%  ' },
    @{ Cell = 'D34'; Value = '5.48' },
    @{ Cell = 'E34'; Value = '  -0.53%  ' },
    @{ Cell = 'E35'; Value = '  +2.80%  ' },
    @{ Cell = 'B36'; Value = 'Kaspa' },
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas' },
    @{ Cell = 'D36'; Value = '0.109' },
    @{ Cell = 'E36'; Value = '  -4.22%  ' },
    @{ Cell = 'B37'; Value = 'VeChain' },
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' },
    @{ Cell = 'D37'; Value = '0.0360' },
    @{ Cell = 'E37'; Value = '  +2.19%  ' },
    @{ Cell = 'D38'; Value = '4.43' },
    @{ Cell = 'E38'; Value = '  +2.47%  ' },
    @{ Cell = 'D39'; Value = '3.39' },
    @{ Cell = 'E39'; Value = '  -5.70%  ' },
    @{ Cell = 'D40'; Value = '0.242' },
    @{ Cell = 'E40'; Value = '  +5.89%  ' },
    @{ Cell = 'D41'; Value = '2.38' },
    @{ Cell = 'E41'; Value = '  +8.55%  ' },
    @{ Cell = 'D42'; Value = '1.39' },
    @{ Cell = 'E42'; Value = '  +20.01%  ' },
    @{ Cell = 'D43'; Value = '12.19' },
    @{ Cell = 'E43'; Value = '  -5.18%  ' },
    @{ Cell = 'D44'; Value = '9.22' },
    @{ Cell = 'E44'; Value = '  +10.38%  ' },
    @{ Cell = 'D45'; Value = '62.51' },
    @{ Cell = 'E45'; Value = '  -2.00%  ' },
    @{ Cell = 'D46'; Value = '5.35' },
    @{ Cell = 'E46'; Value = '  -0.26%  ' },
    @{ Cell = 'E47'; Value = '  +2.78%  ' },
    @{ Cell = 'D48'; Value = '100.40' },
    @{ Cell = 'E48'; Value = '  -0.06%  ' },
    @{ Cell = 'E49'; Value = '  +0.43%  ' },
    @{ Cell = 'D50'; Value = '0.191' },
    @{ Cell = 'E50'; Value = '  +15.92%  ' },
    @{ Cell = 'D51'; Value = '2.543.25' },
    @{ Cell = 'E51'; Value = '  +3.42%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $needsTextFormat = ($u.Cell.Substring(0,1) -eq 'D') -and ($u.Value -match '^[+-]?[0-9]*\.?[0-9]+$')
    if ($needsTextFormat) {
        $range.NumberFormat = '@'
        $range.Value = $u.Value
        $range.ClearFormats()
    } else {
        $range.Value = $u.Value
    }
}
